# Apply "想去人数" (F column) value updates to sheets "展览" (sheet1) and "全部类型" (sheet4)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 885  # F2: 883 -> 885
$ws.Cells.Item(5, 6).Value = 1187  # F5: 1185 -> 1187
$ws.Cells.Item(6, 6).Value = 68  # F6: 67 -> 68
$ws.Cells.Item(7, 6).Value = 4355  # F7: 4346 -> 4355
$ws.Cells.Item(8, 6).Value = 2595  # F8: 2592 -> 2595
$ws.Cells.Item(10, 6).Value = 2508  # F10: 2505 -> 2508
$ws.Cells.Item(14, 6).Value = 1657  # F14: 1655 -> 1657
$ws.Cells.Item(15, 6).Value = 660  # F15: 658 -> 660
$ws.Cells.Item(17, 6).Value = 110  # F17: 109 -> 110
$ws.Cells.Item(18, 6).Value = 321  # F18: 319 -> 321
$ws.Cells.Item(23, 6).Value = 478  # F23: 476 -> 478
$ws.Cells.Item(26, 6).Value = 542  # F26: 540 -> 542
$ws.Cells.Item(29, 6).Value = 80  # F29: 79 -> 80
$ws.Cells.Item(30, 6).Value = 403  # F30: 400 -> 403
$ws.Cells.Item(31, 6).Value = 48  # F31: 47 -> 48
$ws.Cells.Item(32, 6).Value = 1616  # F32: 1615 -> 1616
$ws.Cells.Item(33, 6).Value = 1018  # F33: 1011 -> 1018
$ws.Cells.Item(34, 6).Value = 120  # F34: 118 -> 120
$ws.Cells.Item(35, 6).Value = 19  # F35: 18 -> 19
$ws.Cells.Item(36, 6).Value = 1113  # F36: 1107 -> 1113
$ws.Cells.Item(37, 6).Value = 2043  # F37: 2039 -> 2043
$ws.Cells.Item(38, 6).Value = 261  # F38: 260 -> 261
$ws.Cells.Item(40, 6).Value = 542  # F40: 540 -> 542
$ws.Cells.Item(43, 6).Value = 654  # F43: 653 -> 654
$ws.Cells.Item(44, 6).Value = 1315  # F44: 1312 -> 1315
$ws.Cells.Item(45, 6).Value = 90  # F45: 89 -> 90
$ws.Cells.Item(47, 6).Value = 431  # F47: 428 -> 431

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 885  # F2: 883 -> 885
$ws.Cells.Item(3, 6).Value = 1187  # F3: 1185 -> 1187
$ws.Cells.Item(5, 6).Value = 68  # F5: 67 -> 68
$ws.Cells.Item(6, 6).Value = 4355  # F6: 4346 -> 4355
$ws.Cells.Item(7, 6).Value = 2595  # F7: 2592 -> 2595
$ws.Cells.Item(8, 6).Value = 2508  # F8: 2505 -> 2508
$ws.Cells.Item(9, 6).Value = 1657  # F9: 1655 -> 1657
$ws.Cells.Item(12, 6).Value = 660  # F12: 658 -> 660
$ws.Cells.Item(14, 6).Value = 110  # F14: 109 -> 110
$ws.Cells.Item(15, 6).Value = 321  # F15: 319 -> 321
$ws.Cells.Item(19, 6).Value = 478  # F19: 476 -> 478
$ws.Cells.Item(22, 6).Value = 542  # F22: 540 -> 542
$ws.Cells.Item(28, 6).Value = 80  # F28: 79 -> 80
$ws.Cells.Item(29, 6).Value = 403  # F29: 400 -> 403
$ws.Cells.Item(30, 6).Value = 1616  # F30: 1615 -> 1616
$ws.Cells.Item(31, 6).Value = 1018  # F31: 1011 -> 1018
$ws.Cells.Item(32, 6).Value = 120  # F32: 118 -> 120
$ws.Cells.Item(34, 6).Value = 2044  # F34: 2039 -> 2044
$ws.Cells.Item(35, 6).Value = 261  # F35: 260 -> 261
$ws.Cells.Item(40, 6).Value = 542  # F40: 540 -> 542
$ws.Cells.Item(43, 6).Value = 654  # F43: 653 -> 654
$ws.Cells.Item(44, 6).Value = 1315  # F44: 1312 -> 1315
$ws.Cells.Item(46, 6).Value = 90  # F46: 89 -> 90
$ws.Cells.Item(47, 6).Value = 431  # F47: 428 -> 431

